$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3518616666666667
$ws.Range("H2").Value = 1.055585
$ws.Range("I2").Value = 0.5958054833396739
$ws.Range("J2").Value = 0.5958054833396738
$ws.Range("M2").Value = 28.19948866666667
$ws.Range("N2").Value = 84.598466
$ws.Range("O2").Value = 0.7357427920402423
$ws.Range("P2").Value = 0.7357427920402422
$ws.Range("Q2").Value = 9.922319081401112
$ws.Range("R2").Value = 89.30087173261
$ws.Range("S2").Value = 0.4383595898252177
$ws.Range("T2").Value = 0.4383595898252176
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3518616666666667
$ws.Range("H3").Value = 1.055585
$ws.Range("I3").Value = 0.5958054833396739
$ws.Range("J3").Value = 0.5958054833396738
$ws.Range("O3").Value = 0.2029336910395279
$ws.Range("P3").Value = 0.2029336910395278
$ws.Range("Q3").Value = 2.736789074449445
$ws.Range("R3").Value = 24.631101670045
$ws.Range("S3").Value = 0.1209090058757099
$ws.Range("T3").Value = 0.1209090058757099
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3518616666666667
$ws.Range("H4").Value = 1.055585
$ws.Range("I4").Value = 0.5958054833396739
$ws.Range("J4").Value = 0.5958054833396738
$ws.Range("M4").Value = 2.350402666666667
$ws.Range("N4").Value = 7.051208000000001
$ws.Range("O4").Value = 0.0613235169202299
$ws.Range("P4").Value = 0.06132351692022989
$ws.Range("Q4").Value = 0.8270165996311113
$ws.Range("R4").Value = 7.443149396680001
$ws.Range("S4").Value = 0.03653688763874625
$ws.Range("T4").Value = 0.03653688763874623
$ws.Range("H5").Value = 0.716109
$ws.Range("I5").Value = 0.4041945166603262
$ws.Range("J5").Value = 0.4041945166603262
$ws.Range("M5").Value = 28.19948866666667
$ws.Range("N5").Value = 84.598466
$ws.Range("O5").Value = 0.7357427920402423
$ws.Range("P5").Value = 0.7357427920402422
$ws.Range("Q5").Value = 6.731302543199334
$ws.Range("R5").Value = 60.581722888794
$ws.Range("S5").Value = 0.2973832022150247
$ws.Range("T5").Value = 0.2973832022150246
$ws.Range("H6").Value = 0.716109
$ws.Range("I6").Value = 0.4041945166603262
$ws.Range("J6").Value = 0.4041945166603262
$ws.Range("O6").Value = 0.2029336910395279
$ws.Range("P6").Value = 0.2029336910395278
$ws.Range("Q6").Value = 1.856638060710333
$ws.Range("S6").Value = 0.08202468516381795
$ws.Range("T6").Value = 0.08202468516381792
$ws.Range("H7").Value = 0.716109
$ws.Range("I7").Value = 0.4041945166603262
$ws.Range("J7").Value = 0.4041945166603262
$ws.Range("M7").Value = 2.350402666666667
$ws.Range("N7").Value = 7.051208000000001
$ws.Range("O7").Value = 0.0613235169202299
$ws.Range("P7").Value = 0.06132351692022989
$ws.Range("Q7").Value = 0.5610481677413335
$ws.Range("S7").Value = 0.02478662928148366
$ws.Range("T7").Value = 0.02478662928148365
